# Add a new worksheet "tc074" right after "tc051" (before "tc007"),
# populate it with the Epic/Feature/Requirement data-provider row pair,
# and leave it as the active/selected sheet (matching the author's
# "Added TC-74 / TC-75 / TC-76 data provider" commit).

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("tc051")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "tc074"

$newSheet.Range("A1").Value = "Epic"
$newSheet.Range("B1").Value = "Feature"
$newSheet.Range("C1").Value = "Requirement"
$newSheet.Range("A2").Value = "Epic Mohit"
$newSheet.Range("B2").Value = "Mohit Feature"
$newSheet.Range("C2").Value = "RQ-489"

$newSheet.Range("B7").Select()
$newSheet.Activate()
